$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.134.06'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '2.980.68'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'595.82"
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("D6").Value = "'143.00"
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D9").Value = '2.980.85'
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").Value = "'5.99"
$ws.Range("E11").Value = '  +5.07%  '
$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = '  +2.70%  '
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").Value = "'34.11"
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '3.469.71'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").Value = '61.212.24'
$ws.Range("E17").Value = '  -1.66%  '
$ws.Range("D18").Value = "'6.85"
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").Value = '2.979.62'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = "'449.19"
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("D21").Value = "'14.00"
$ws.Range("E21").Value = '  +1.31%  '
$ws.Range("D22").Value = "'0.681"
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").Value = "'7.28"
$ws.Range("E23").Value = '  -0.81%  '
$ws.Range("E24").Value = '  +2.40%  '
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = '  -3.43%  '
$ws.Range("D26").Value = "'10.47"
$ws.Range("E26").Value = '  +3.24%  '
$ws.Range("D27").Value = "'11.92"
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  +2.52%  '
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").Value = "'7.14"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").Value = "'2.05"
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("D33").Value = "'27.25"
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").Value = '0.0₃0809'
$ws.Range("E35").Value = '  +3.31%  '
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = "'5.77"
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("D38").Value = "'50.02"
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").Value = "'2.04"
$ws.Range("E39").Value = '  -2.90%  '
$ws.Range("D40").Value = "'8.96"
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("E41").Value = '  +6.02%  '
$ws.Range("D42").Value = "'2.83"
$ws.Range("E42").Value = '  -2.92%  '
$ws.Range("D43").Value = "'387.88"
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("D44").Value = "'0.268"
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("D45").Value = "'0.0348"
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").Value = "'38.60"
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("D47").Value = '2.698.85'
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("D51").Value = "'2.13"
$ws.Range("E51").Value = '  -0.13%  '
